$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 5.441829000000001
$ws.Cells.Item(2, 8).Value = 16.325487
$ws.Cells.Item(2, 9).Value = 0.5729403216841985
$ws.Cells.Item(2, 10).Value = 0.5729403216841985
$ws.Cells.Item(2, 13).Value = 181.3526613333333
$ws.Cells.Item(2, 14).Value = 544.057984
$ws.Cells.Item(2, 15).Value = 0.9845849379007657
$ws.Cells.Item(2, 16).Value = 0.984584937900766
$ws.Cells.Item(2, 17).Value = 986.8901716709123
$ws.Cells.Item(2, 18).Value = 8882.011545038209
$ws.Cells.Item(2, 19).Value = 0.5641084110462813
$ws.Cells.Item(2, 20).Value = 0.5641084110462815
# Row 3
$ws.Cells.Item(3, 7).Value = 5.441829000000001
$ws.Cells.Item(3, 8).Value = 16.325487
$ws.Cells.Item(3, 9).Value = 0.5729403216841985
$ws.Cells.Item(3, 10).Value = 0.5729403216841985
$ws.Cells.Item(3, 15).Value = 0.003686045149950483
$ws.Cells.Item(3, 16).Value = 0.003686045149950484
$ws.Cells.Item(3, 17).Value = 3.694675381260001
$ws.Cells.Item(3, 18).Value = 33.25207843134
$ws.Cells.Item(3, 19).Value = 0.00211188389395511
$ws.Cells.Item(3, 20).Value = 0.00211188389395511
# Row 4
$ws.Cells.Item(4, 7).Value = 5.441829000000001
$ws.Cells.Item(4, 8).Value = 16.325487
$ws.Cells.Item(4, 9).Value = 0.5729403216841985
$ws.Cells.Item(4, 10).Value = 0.5729403216841985
$ws.Cells.Item(4, 13).Value = 0.6398506666666667
$ws.Cells.Item(4, 14).Value = 1.919552
$ws.Cells.Item(4, 15).Value = 0.003473824559694892
$ws.Cells.Item(4, 16).Value = 0.003473824559694892
$ws.Cells.Item(4, 17).Value = 3.481957913536001
$ws.Cells.Item(4, 18).Value = 31.337621221824
$ws.Cells.Item(4, 19).Value = 0.001990294160706061
$ws.Cells.Item(4, 20).Value = 0.001990294160706061
# Row 5
$ws.Cells.Item(5, 7).Value = 5.441829000000001
$ws.Cells.Item(5, 8).Value = 16.325487
$ws.Cells.Item(5, 9).Value = 0.5729403216841985
$ws.Cells.Item(5, 10).Value = 0.5729403216841985
$ws.Cells.Item(5, 13).Value = 1.520540333333333
$ws.Cells.Item(5, 14).Value = 4.561621
$ws.Cells.Item(5, 15).Value = 0.008255192389588805
$ws.Cells.Item(5, 16).Value = 0.008255192389588807
$ws.Cells.Item(5, 17).Value = 8.274520481603
$ws.Cells.Item(5, 18).Value = 74.47068433442701
$ws.Cells.Item(5, 19).Value = 0.004729732583255958
$ws.Cells.Item(5, 20).Value = 0.004729732583255959
# Row 6
$ws.Cells.Item(6, 9).Value = 0.2716201486343598
$ws.Cells.Item(6, 10).Value = 0.2716201486343598
$ws.Cells.Item(6, 13).Value = 181.3526613333333
$ws.Cells.Item(6, 14).Value = 544.057984
$ws.Cells.Item(6, 15).Value = 0.9845849379007657
$ws.Cells.Item(6, 16).Value = 0.984584937900766
$ws.Cells.Item(6, 17).Value = 467.865927688704
$ws.Cells.Item(6, 18).Value = 4210.793349198336
$ws.Cells.Item(6, 19).Value = 0.2674331071757579
$ws.Cells.Item(6, 20).Value = 0.267433107175758
# Row 7
$ws.Cells.Item(7, 9).Value = 0.2716201486343598
$ws.Cells.Item(7, 10).Value = 0.2716201486343598
$ws.Cells.Item(7, 15).Value = 0.003686045149950483
$ws.Cells.Item(7, 16).Value = 0.003686045149950484
$ws.Cells.Item(7, 19).Value = 0.001001204131502511
$ws.Cells.Item(7, 20).Value = 0.001001204131502511
# Row 8
$ws.Cells.Item(8, 9).Value = 0.2716201486343598
$ws.Cells.Item(8, 10).Value = 0.2716201486343598
$ws.Cells.Item(8, 13).Value = 0.6398506666666667
$ws.Cells.Item(8, 14).Value = 1.919552
$ws.Cells.Item(8, 15).Value = 0.003473824559694892
$ws.Cells.Item(8, 16).Value = 0.003473824559694892
$ws.Cells.Item(8, 17).Value = 1.650730259712
$ws.Cells.Item(8, 18).Value = 14.856572337408
$ws.Cells.Item(8, 19).Value = 0.000943560743234016
$ws.Cells.Item(8, 20).Value = 0.0009435607432340161
# Row 9
$ws.Cells.Item(9, 9).Value = 0.2716201486343598
$ws.Cells.Item(9, 10).Value = 0.2716201486343598
$ws.Cells.Item(9, 13).Value = 1.520540333333333
$ws.Cells.Item(9, 14).Value = 4.561621
$ws.Cells.Item(9, 15).Value = 0.008255192389588805
$ws.Cells.Item(9, 16).Value = 0.008255192389588807
$ws.Cells.Item(9, 17).Value = 3.922793348675999
$ws.Cells.Item(9, 18).Value = 35.305140138084
$ws.Cells.Item(9, 19).Value = 0.002242276583865347
$ws.Cells.Item(9, 20).Value = 0.002242276583865348
# Row 10
$ws.Cells.Item(10, 7).Value = 1.476376
$ws.Cells.Item(10, 8).Value = 4.429128
$ws.Cells.Item(10, 9).Value = 0.1554395296814417
$ws.Cells.Item(10, 10).Value = 0.1554395296814417
$ws.Cells.Item(10, 13).Value = 181.3526613333333
$ws.Cells.Item(10, 14).Value = 544.057984
$ws.Cells.Item(10, 15).Value = 0.9845849379007657
$ws.Cells.Item(10, 16).Value = 0.984584937900766
$ws.Cells.Item(10, 17).Value = 267.7447167286614
$ws.Cells.Item(10, 18).Value = 2409.702450557952
$ws.Cells.Item(10, 19).Value = 0.1530434196787265
$ws.Cells.Item(10, 20).Value = 0.1530434196787265
# Row 11
$ws.Cells.Item(11, 7).Value = 1.476376
$ws.Cells.Item(11, 8).Value = 4.429128
$ws.Cells.Item(11, 9).Value = 0.1554395296814417
$ws.Cells.Item(11, 10).Value = 0.1554395296814417
$ws.Cells.Item(11, 15).Value = 0.003686045149950483
$ws.Cells.Item(11, 16).Value = 0.003686045149950484
$ws.Cells.Item(11, 17).Value = 1.00237072144
$ws.Cells.Item(11, 18).Value = 9.021336492960002
$ws.Cells.Item(11, 19).Value = 0.0005729571244928622
$ws.Cells.Item(11, 20).Value = 0.0005729571244928625
# Row 12
$ws.Cells.Item(12, 7).Value = 1.476376
$ws.Cells.Item(12, 8).Value = 4.429128
$ws.Cells.Item(12, 9).Value = 0.1554395296814417
$ws.Cells.Item(12, 10).Value = 0.1554395296814417
$ws.Cells.Item(12, 13).Value = 0.6398506666666667
$ws.Cells.Item(12, 14).Value = 1.919552
$ws.Cells.Item(12, 15).Value = 0.003473824559694892
$ws.Cells.Item(12, 16).Value = 0.003473824559694892
$ws.Cells.Item(12, 17).Value = 0.9446601678506668
$ws.Cells.Item(12, 18).Value = 8.501941510656
$ws.Cells.Item(12, 19).Value = 0.0005399696557548152
$ws.Cells.Item(12, 20).Value = 0.0005399696557548154
# Row 13
$ws.Cells.Item(13, 7).Value = 1.476376
$ws.Cells.Item(13, 8).Value = 4.429128
$ws.Cells.Item(13, 9).Value = 0.1554395296814417
$ws.Cells.Item(13, 10).Value = 0.1554395296814417
$ws.Cells.Item(13, 13).Value = 1.520540333333333
$ws.Cells.Item(13, 14).Value = 4.561621
$ws.Cells.Item(13, 15).Value = 0.008255192389588805
$ws.Cells.Item(13, 16).Value = 0.008255192389588807
$ws.Cells.Item(13, 17).Value = 2.244889255165333
$ws.Cells.Item(13, 18).Value = 20.204003296488
$ws.Cells.Item(13, 19).Value = 0.001283183222467501
$ws.Cells.Item(13, 20).Value = 0.001283183222467501
